$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the misspelling "Area, Sq.fit" -> "Area, Sq.ft" in the header cell E1
$ws.Range("E1").Value = "Area, Sq.ft"

# Update the active cell selection to match the saved view state
$ws.Range("H8").Select() | Out-Null
